$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A212:A249").EntireRow.Delete()
